$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("experiments")

# Row 386
$ws.Cells.Item(386, 14).Value = 1.0
$ws.Cells.Item(386, 15).Value = 387.0
$ws.Cells.Item(386, 16).Value = 0.04396420463444639
$ws.Cells.Item(386, 17).Value = 17.238508626221513
$ws.Cells.Item(386, 18).Value = 77.19298245614036
$ws.Cells.Item(386, 18).Style = "Normal"
$ws.Cells.Item(386, 19).Value = 73.6842105263158
$ws.Cells.Item(386, 19).Style = "Normal"
$ws.Cells.Item(386, 20).Value = 'T3_Supination and pronation_Sammon_FS-IF-IA_06-Apr-2017_MI_3000_LR_0.5'
$ws.Cells.Item(386, 21).Value = 'ResultsMar2017-Sammon_MI_3000_LR_0.5.xlsx'

# Row 387
$ws.Cells.Item(387, 14).Value = 1.0
$ws.Cells.Item(387, 15).Value = 561.0
$ws.Cells.Item(387, 16).Value = 0.025518751723765307
$ws.Cells.Item(387, 17).Value = 7.450914789672176
$ws.Cells.Item(387, 18).Value = 68.42105263157895
$ws.Cells.Item(387, 18).Style = "Normal"
$ws.Cells.Item(387, 19).Value = 63.15789473684212
$ws.Cells.Item(387, 19).Style = "Normal"
$ws.Cells.Item(387, 20).Value = 'T4_Rest_Sammon_FS_06-Apr-2017_MI_3000_LR_0.5'
$ws.Cells.Item(387, 21).Value = 'ResultsMar2017-Sammon_MI_3000_LR_0.5.xlsx'

# Row 388
$ws.Cells.Item(388, 14).Value = 1.0
$ws.Cells.Item(388, 15).Value = 467.0
$ws.Cells.Item(388, 16).Value = 0.04743784469949855
$ws.Cells.Item(388, 17).Value = 17.747292084363693
$ws.Cells.Item(388, 18).Value = 70.17543859649123
$ws.Cells.Item(388, 18).Style = "Normal"
$ws.Cells.Item(388, 19).Value = 78.94736842105263
$ws.Cells.Item(388, 19).Style = "Normal"
$ws.Cells.Item(388, 20).Value = 'T4_Rest_Sammon_IF_06-Apr-2017_MI_3000_LR_0.5'
$ws.Cells.Item(388, 21).Value = 'ResultsMar2017-Sammon_MI_3000_LR_0.5.xlsx'

# Row 389
$ws.Cells.Item(389, 14).Value = 1.0
$ws.Cells.Item(389, 15).Value = 348.0
$ws.Cells.Item(389, 16).Value = 0.02871769084791272
$ws.Cells.Item(389, 17).Value = 29.386682582976263
$ws.Cells.Item(389, 18).Value = 67.83625730994152
$ws.Cells.Item(389, 18).Style = "Normal"
$ws.Cells.Item(389, 19).Value = 78.94736842105263
$ws.Cells.Item(389, 19).Style = "Normal"
$ws.Cells.Item(389, 20).Value = 'T4_Rest_Sammon_IA_06-Apr-2017_MI_3000_LR_0.5'
$ws.Cells.Item(389, 21).Value = 'ResultsMar2017-Sammon_MI_3000_LR_0.5.xlsx'

# Row 390
$ws.Cells.Item(390, 14).Value = 1.0
$ws.Cells.Item(390, 15).Value = 3000.0
$ws.Cells.Item(390, 16).Value = 0.040678362939114684
$ws.Cells.Item(390, 17).Value = 21.370258284066885
$ws.Cells.Item(390, 18).Value = 77.77777777777777
$ws.Cells.Item(390, 18).Style = "Normal"
$ws.Cells.Item(390, 19).Value = 78.94736842105263
$ws.Cells.Item(390, 19).Style = "Normal"
$ws.Cells.Item(390, 20).Value = 'T4_Rest_Sammon_FS-IF_06-Apr-2017_MI_3000_LR_0.5'
$ws.Cells.Item(390, 21).Value = 'ResultsMar2017-Sammon_MI_3000_LR_0.5.xlsx'

# Row 391
$ws.Cells.Item(391, 14).Value = 1.0
$ws.Cells.Item(391, 15).Value = 495.0
$ws.Cells.Item(391, 16).Value = 0.028080396363575802
$ws.Cells.Item(391, 17).Value = 17.439436136485227
$ws.Cells.Item(391, 18).Value = 70.76023391812865
$ws.Cells.Item(391, 18).Style = "Normal"
$ws.Cells.Item(391, 19).Value = 73.6842105263158
$ws.Cells.Item(391, 19).Style = "Normal"
$ws.Cells.Item(391, 20).Value = 'T4_Rest_Sammon_FS-IA_06-Apr-2017_MI_3000_LR_0.5'
$ws.Cells.Item(391, 21).Value = 'ResultsMar2017-Sammon_MI_3000_LR_0.5.xlsx'

# Row 392
$ws.Cells.Item(392, 14).Value = 1.0
$ws.Cells.Item(392, 15).Value = 1021.0
$ws.Cells.Item(392, 16).Value = 0.043332922454232536
$ws.Cells.Item(392, 17).Value = 27.819165346437043
$ws.Cells.Item(392, 18).Value = 75.43859649122807
$ws.Cells.Item(392, 18).Style = "Normal"
$ws.Cells.Item(392, 19).Value = 78.94736842105263
$ws.Cells.Item(392, 19).Style = "Normal"
$ws.Cells.Item(392, 20).Value = 'T4_Rest_Sammon_IF-IA_06-Apr-2017_MI_3000_LR_0.5'
$ws.Cells.Item(392, 21).Value = 'ResultsMar2017-Sammon_MI_3000_LR_0.5.xlsx'

# Row 393
$ws.Cells.Item(393, 14).Value = 1.0
$ws.Cells.Item(393, 15).Value = 826.0
$ws.Cells.Item(393, 16).Value = 0.040010344537007976
$ws.Cells.Item(393, 17).Value = 23.456215257883798
$ws.Cells.Item(393, 18).Value = 75.43859649122807
$ws.Cells.Item(393, 18).Style = "Normal"
$ws.Cells.Item(393, 19).Value = 78.94736842105263
$ws.Cells.Item(393, 19).Style = "Normal"
$ws.Cells.Item(393, 20).Value = 'T4_Rest_Sammon_FS-IF-IA_06-Apr-2017_MI_3000_LR_0.5'
$ws.Cells.Item(393, 21).Value = 'ResultsMar2017-Sammon_MI_3000_LR_0.5.xlsx'

# Row 394
$ws.Cells.Item(394, 14).Value = 1.0
$ws.Cells.Item(394, 15).Value = 386.0
$ws.Cells.Item(394, 16).Value = 0.033276475332387735
$ws.Cells.Item(394, 17).Value = 7.492143424751974
$ws.Cells.Item(394, 18).Value = 85.38011695906432
$ws.Cells.Item(394, 18).Style = "Normal"
$ws.Cells.Item(394, 19).Value = 68.42105263157896
$ws.Cells.Item(394, 19).Style = "Normal"
$ws.Cells.Item(394, 20).Value = 'T1_Finger taps_Sammon_FS_06-Apr-2017_MI_4000_LR_0.5'
$ws.Cells.Item(394, 21).Value = 'ResultsMar2017-Sammon_MI_4000_LR_0.5.xlsx'

# Row 395
$ws.Cells.Item(395, 14).Value = 1.0
$ws.Cells.Item(395, 15).Value = 425.0
$ws.Cells.Item(395, 16).Value = 0.05297079079312436
$ws.Cells.Item(395, 17).Value = 22.26747779750949
$ws.Cells.Item(395, 18).Value = 73.09941520467837
$ws.Cells.Item(395, 18).Style = "Normal"
$ws.Cells.Item(395, 19).Value = 52.63157894736844
$ws.Cells.Item(395, 19).Style = "Normal"
$ws.Cells.Item(395, 20).Value = 'T1_Finger taps_Sammon_IF_06-Apr-2017_MI_4000_LR_0.5'
$ws.Cells.Item(395, 21).Value = 'ResultsMar2017-Sammon_MI_4000_LR_0.5.xlsx'

# Row 396
$ws.Cells.Item(396, 14).Value = 1.0
$ws.Cells.Item(396, 15).Value = 543.0
$ws.Cells.Item(396, 16).Value = 0.03474891418136164
$ws.Cells.Item(396, 17).Value = 6.934086693764326
$ws.Cells.Item(396, 18).Value = 78.94736842105263
$ws.Cells.Item(396, 18).Style = "Normal"
$ws.Cells.Item(396, 19).Value = 68.42105263157896
$ws.Cells.Item(396, 19).Style = "Normal"
$ws.Cells.Item(396, 20).Value = 'T1_Finger taps_Sammon_IA_06-Apr-2017_MI_4000_LR_0.5'
$ws.Cells.Item(396, 21).Value = 'ResultsMar2017-Sammon_MI_4000_LR_0.5.xlsx'

# Row 397
$ws.Cells.Item(397, 14).Value = 1.0
$ws.Cells.Item(397, 15).Value = 1461.0
$ws.Cells.Item(397, 16).Value = 0.04757995748830339
$ws.Cells.Item(397, 17).Value = 38.79562766488739
$ws.Cells.Item(397, 18).Value = 80.11695906432749
$ws.Cells.Item(397, 18).Style = "Normal"
$ws.Cells.Item(397, 19).Value = 68.42105263157896
$ws.Cells.Item(397, 19).Style = "Normal"
$ws.Cells.Item(397, 20).Value = 'T1_Finger taps_Sammon_FS-IF_06-Apr-2017_MI_4000_LR_0.5'
$ws.Cells.Item(397, 21).Value = 'ResultsMar2017-Sammon_MI_4000_LR_0.5.xlsx'

# Row 398
$ws.Cells.Item(398, 14).Value = 1.0
$ws.Cells.Item(398, 15).Value = 668.0
$ws.Cells.Item(398, 16).Value = 0.035779420904758184
$ws.Cells.Item(398, 17).Value = 10.651671977948329
$ws.Cells.Item(398, 18).Value = 81.28654970760233
$ws.Cells.Item(398, 18).Style = "Normal"
$ws.Cells.Item(398, 19).Value = 73.6842105263158
$ws.Cells.Item(398, 19).Style = "Normal"
$ws.Cells.Item(398, 20).Value = 'T1_Finger taps_Sammon_FS-IA_06-Apr-2017_MI_4000_LR_0.5'
$ws.Cells.Item(398, 21).Value = 'ResultsMar2017-Sammon_MI_4000_LR_0.5.xlsx'

# Row 399
$ws.Cells.Item(399, 14).Value = 1.0
$ws.Cells.Item(399, 15).Value = 1902.0
$ws.Cells.Item(399, 16).Value = 0.04852974848121782
$ws.Cells.Item(399, 17).Value = 17.63774121897484
$ws.Cells.Item(399, 18).Value = 80.11695906432749
$ws.Cells.Item(399, 18).Style = "Normal"
$ws.Cells.Item(399, 19).Value = 63.15789473684212
$ws.Cells.Item(399, 19).Style = "Normal"
$ws.Cells.Item(399, 20).Value = 'T1_Finger taps_Sammon_IF-IA_06-Apr-2017_MI_4000_LR_0.5'
$ws.Cells.Item(399, 21).Value = 'ResultsMar2017-Sammon_MI_4000_LR_0.5.xlsx'

# Row 400
$ws.Cells.Item(400, 14).Value = 1.0
$ws.Cells.Item(400, 15).Value = 567.0
$ws.Cells.Item(400, 16).Value = 0.04470052095758281
$ws.Cells.Item(400, 17).Value = 34.20937973731283
$ws.Cells.Item(400, 18).Value = 81.28654970760233
$ws.Cells.Item(400, 18).Style = "Normal"
$ws.Cells.Item(400, 19).Value = 84.21052631578948
$ws.Cells.Item(400, 19).Style = "Normal"
$ws.Cells.Item(400, 20).Value = 'T1_Finger taps_Sammon_FS-IF-IA_06-Apr-2017_MI_4000_LR_0.5'
$ws.Cells.Item(400, 21).Value = 'ResultsMar2017-Sammon_MI_4000_LR_0.5.xlsx'

# Row 401
$ws.Cells.Item(401, 14).Value = 1.0
$ws.Cells.Item(401, 15).Value = 359.0
$ws.Cells.Item(401, 16).Value = 0.03623810463778771
$ws.Cells.Item(401, 17).Value = 3.6522599069320476
$ws.Cells.Item(401, 18).Value = 80.11695906432749
$ws.Cells.Item(401, 18).Style = "Normal"
$ws.Cells.Item(401, 19).Value = 73.6842105263158
$ws.Cells.Item(401, 19).Style = "Normal"
$ws.Cells.Item(401, 20).Value = 'T2_Finger to nose_Sammon_FS_06-Apr-2017_MI_4000_LR_0.5'
$ws.Cells.Item(401, 21).Value = 'ResultsMar2017-Sammon_MI_4000_LR_0.5.xlsx'

# Row 402
$ws.Cells.Item(402, 14).Value = 1.0
$ws.Cells.Item(402, 15).Value = 577.0
$ws.Cells.Item(402, 16).Value = 0.05189347750460582
$ws.Cells.Item(402, 17).Value = 17.2051003349763
$ws.Cells.Item(402, 18).Value = 77.19298245614036
$ws.Cells.Item(402, 18).Style = "Normal"
$ws.Cells.Item(402, 19).Value = 73.6842105263158
$ws.Cells.Item(402, 19).Style = "Normal"
$ws.Cells.Item(402, 20).Value = 'T2_Finger to nose_Sammon_IF_06-Apr-2017_MI_4000_LR_0.5'
$ws.Cells.Item(402, 21).Value = 'ResultsMar2017-Sammon_MI_4000_LR_0.5.xlsx'

# Row 403
$ws.Cells.Item(403, 14).Value = 1.0
$ws.Cells.Item(403, 15).Value = 399.0
$ws.Cells.Item(403, 16).Value = 0.035121183484197536
$ws.Cells.Item(403, 17).Value = 7.0167752252638085
$ws.Cells.Item(403, 18).Value = 75.43859649122807
$ws.Cells.Item(403, 18).Style = "Normal"
$ws.Cells.Item(403, 19).Value = 68.42105263157896
$ws.Cells.Item(403, 19).Style = "Normal"
$ws.Cells.Item(403, 20).Value = 'T2_Finger to nose_Sammon_IA_06-Apr-2017_MI_4000_LR_0.5'
$ws.Cells.Item(403, 21).Value = 'ResultsMar2017-Sammon_MI_4000_LR_0.5.xlsx'

# Row 404
$ws.Cells.Item(404, 14).Value = 1.0
$ws.Cells.Item(404, 15).Value = 351.0
$ws.Cells.Item(404, 16).Value = 0.04867543591896745
$ws.Cells.Item(404, 17).Value = 21.371173527575458
$ws.Cells.Item(404, 18).Value = 85.96491228070175
$ws.Cells.Item(404, 18).Style = "Normal"
$ws.Cells.Item(404, 19).Value = 84.21052631578948
$ws.Cells.Item(404, 19).Style = "Normal"
$ws.Cells.Item(404, 20).Value = 'T2_Finger to nose_Sammon_FS-IF_06-Apr-2017_MI_4000_LR_0.5'
$ws.Cells.Item(404, 21).Value = 'ResultsMar2017-Sammon_MI_4000_LR_0.5.xlsx'

# Row 405
$ws.Cells.Item(405, 14).Value = 1.0
$ws.Cells.Item(405, 15).Value = 491.0
$ws.Cells.Item(405, 16).Value = 0.03865048214386203
$ws.Cells.Item(405, 17).Value = 7.9605961691835
$ws.Cells.Item(405, 18).Value = 70.17543859649123
$ws.Cells.Item(405, 18).Style = "Normal"
$ws.Cells.Item(405, 19).Value = 84.21052631578948
$ws.Cells.Item(405, 19).Style = "Normal"
$ws.Cells.Item(405, 20).Value = 'T2_Finger to nose_Sammon_FS-IA_06-Apr-2017_MI_4000_LR_0.5'
$ws.Cells.Item(405, 21).Value = 'ResultsMar2017-Sammon_MI_4000_LR_0.5.xlsx'

# Row 406
$ws.Cells.Item(406, 14).Value = 1.0
$ws.Cells.Item(406, 15).Value = 509.0
$ws.Cells.Item(406, 16).Value = 0.048787228405706685
$ws.Cells.Item(406, 17).Value = 21.24894109137171
$ws.Cells.Item(406, 18).Value = 87.71929824561404
$ws.Cells.Item(406, 18).Style = "Normal"
$ws.Cells.Item(406, 19).Value = 68.42105263157896
$ws.Cells.Item(406, 19).Style = "Normal"
$ws.Cells.Item(406, 20).Value = 'T2_Finger to nose_Sammon_IF-IA_06-Apr-2017_MI_4000_LR_0.5'
$ws.Cells.Item(406, 21).Value = 'ResultsMar2017-Sammon_MI_4000_LR_0.5.xlsx'

# Row 407
$ws.Cells.Item(407, 15).Value = 371.0
$ws.Cells.Item(407, 16).Value = 0.04789382389439819
$ws.Cells.Item(407, 20).Value = 'T2_Finger to nose_Sammon_FS-IF-IA_06-Apr-2017_MI_4000_LR_0.5'
